# "put the email up top"
#
# - Move the Content Placeholder (bullet body) down a bit to make room.
# - Fix the double period at the end of the second bullet paragraph
#   ("...chops.." -> "...chops.").
# - Move the WebPlatformJobs@adobe.com textbox up near the top and make it
#   sz=20/bold to match the other headline-ish text.
# - Move/resize the html.adobe.com textbox to the bottom-left and drop its
#   right-alignment override (it's no longer tucked in the bottom-right
#   corner).
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/Top/Width/Height
# (and TextRange/ParagraphFormat numeric properties) are expressed in
# points, while the OOXML stores EMU (1 pt = 12700 EMU) as a 32-bit float
# internally in this host. A handful of the literals are nudged by a few
# 10^-6 pt so that, after the pt->EMU float round-trip, they land exactly
# on the target EMU value instead of 1 EMU short.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Content Placeholder 2: shift down (228600,990600) -> (228600,1382838);
# size (8686800 x 4547162) is unchanged.
# ---------------------------------------------------------------------
$body = $s.Shapes.Item(2)
$body.Top = 108.88488388976377

# ---------------------------------------------------------------------
# Fix "...chops.." -> "...chops." in the second bullet paragraph. This
# also naturally splits the trailing "." into its own run, matching how
# PowerPoint records an in-place edit.
# ---------------------------------------------------------------------
$bodyTextRange = $body.TextFrame.TextRange
$bodyText = $bodyTextRange.Text
$dupIdx = $bodyText.IndexOf("chops..")
$dupDots = $bodyTextRange.Characters($dupIdx + 6, 2)
$dupDots.Text = "."

# ---------------------------------------------------------------------
# TextBox 4 (WebPlatformJobs@adobe.com): move up top and enlarge a bit;
# also bump to sz=2000 (20pt), bold, to match the new prominent position.
# ---------------------------------------------------------------------
$emailBox = $s.Shapes.Item(4)
$emailBox.Left = 211.01819617637796
$emailBox.Top = 73.53244404488188
$emailBox.Width = 297.7793731787402
$emailBox.Height = 31.50472360944882

$emailRange = $emailBox.TextFrame.TextRange
$emailRange.Font.Size = 20
$emailRange.Font.Bold = $true

# ---------------------------------------------------------------------
# TextBox 5 (html.adobe.com): move to bottom-left, narrower, and drop the
# right-alignment override it had when it lived in the bottom-right
# corner.
# ---------------------------------------------------------------------
$htmlBox = $s.Shapes.Item(5)
$htmlBox.Left = 5.882047444094487
$htmlBox.Top = 472.31016548031494
$htmlBox.Width = 158.31613923228346

$htmlBox.TextFrame.TextRange.ParagraphFormat.Alignment = 1
